# Apply updates to the "Instructores" worksheet:
#  - resize columns A, B, C
#  - replace the data in row 2
#  - append new instructor rows 3-16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust column widths (ColumnWidth round-trips through Excel's
# internal padding, which adds ~0.8333 chars back on save, so the
# assigned values are offset to land exactly on the target widths).
$ws.Columns.Item(1).ColumnWidth = 13.1666666666667
$ws.Columns.Item(2).ColumnWidth = 14.1666666666667
$ws.Columns.Item(3).ColumnWidth = 36.1666666666667

# Data rows (nombres, apellidos, correo_institucional, numero_celular, numero_cedula)
$rows = @(
    @("Jesus Andres", "Silva Plazas", "jsapp@sena.edu.co", "3148777918", "1077848000"),
    @("Sebastian", "Ramirez Rojas", "sebastian_ramirezro@soy.sena.edu.co", "3133516648", "1077848001"),
    @("Jose", "Ramirez", "joseramirez@sena.edu.co", "3118777918", "1077488000"),
    @("Maria", "Rodriguez", "joseramirez@sena.edu.co", "3157894562", "1082654789"),
    @("Carlos", "Garcia", "joseramirez@sena.edu.co", "3209876543", "1098765432"),
    @("Ana", "Martinez", "joseramirez@sena.edu.co", "3112345678", "1045678912"),
    @("Juan", "Lopez", "joseramirez@sena.edu.co", "3145678901", "1034567890"),
    @("Patricia", "Gomez", "joseramirez@sena.edu.co", "3187654321", "1067891234"),
    @("Diego", "Hernandez", "joseramirez@sena.edu.co", "3198765432", "1056789123"),
    @("Laura", "Sanchez", "joseramirez@sena.edu.co", "3167891234", "1023456789"),
    @("Miguel", "Torres", "joseramirez@sena.edu.co", "3178912345", "1012345678"),
    @("Carmen", "Perez", "joseramirez@sena.edu.co", "3123456789", "1089123456"),
    @("Ricardo", "Diaz", "joseramirez@sena.edu.co", "3134567890", "1078912345"),
    @("Valeria", "Jimenez", "joseramirez@sena.edu.co", "3190123456", "1001234567"),
    @("Carolina", "Reyes", "joseramirez@sena.edu.co", "3212345678", "1087654321")
)

# Force the phone/ID columns to remain plain text so long numeric
# strings aren't reinterpreted as numbers (they keep leading context
# like '3...' / '1...' verbatim, same as the source inline strings).
$ws.Range("D2:E16").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
